# Fix the "B to B Connector" label on the Block Diagram slide:
#   - correct the wording to "Board to Board Connector"
#   - shrink the font from 10pt to 8.5pt to fit the longer text
#   - nudge/resize the (rotated) text box to its new position/size

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(118)   # "テキスト ボックス 294" - the B to B Connector label

# Correct the label text and reduce its font size to suit the new width.
# (Do this before touching Left/Top/Width/Height: the textbox has
# spAutoFit, so its height is recomputed from the text - setting the
# size first means our explicit Height below is the one that sticks.)
$tr = $sh.TextFrame.TextRange
$tr.Text = "Board to Board Connector"
$tr.Font.Size = 8.5

# Reposition / resize the rotated textbox (values are the unrotated
# bounding box, in points, matching the underlying a:off / a:ext).
$sh.Left   = 71.70465
$sh.Top    = 238.4763
$sh.Width  = 149.8827
$sh.Height = 17.57
